$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Jimy" (previously in A5) becomes "Jimy Carter" and moves to B1
$ws.Range("B1").Value = "Jimy Carter"

# Remove the now-empty trailing row (A5), shrinking the used range to A1:B4
$ws.Range("A5").ClearContents()

# Update the active selection to match the new layout
$ws.Range("B5").Select() | Out-Null
